$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.483.04'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '2.629.18'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '112.27'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = '323.82'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '0.525'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('D10').Value = '39.59'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('D11').Value = '19.75'
$ws.Range('E11').Value = '  -4.16%  '
$ws.Range('D12').Value = '0.0811'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').Value = '7.29'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '3.037.42'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '2.614.09'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '0.851'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '49.426.11'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '12.89'
$ws.Range('E19').Value = '  -3.41%  '
$ws.Range('E20').Value = '  -4.07%  '
$ws.Range('D21').Value = '6.68'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '0.0₃0946'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').Value = '269.41'
$ws.Range('E23').Value = '  -4.50%  '
$ws.Range('D24').Value = '68.94'
$ws.Range('E24').Value = '  -5.71%  '
$ws.Range('E25').Value = '  -2.64%  '
$ws.Range('D26').Value = '26.24'
$ws.Range('E26').Value = '  -2.37%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('D31').Value = '34.69'
$ws.Range('E31').Value = '  -5.46%  '
$ws.Range('D32').Value = '49.42'
$ws.Range('D33').Value = '5.50'
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').Value = '0.0813'
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('D36').Value = '18.81'
$ws.Range('E36').Value = '  -3.85%  '
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').Value = '128.62'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').Value = '0.111'
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('D42').Value = '22.19'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('E43').Value = '  +3.82%  '
$ws.Range('E44').Value = '  -3.73%  '
$ws.Range('D45').Value = '2.053.65'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('E46').Value = '  -5.16%  '
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  +4.96%  '
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '58.71'
$ws.Range('E50').Value = '  +1.14%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '5.19'
$ws.Range('E51').Value = '  -3.95%  '
